$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates per the diff. D-column "Price" cells are forced to Text format
# (NumberFormat "@") before assignment so numeric-looking strings (e.g. "602.20",
# "1.00", "66.775.03") stay exactly as authored instead of Excel auto-converting
# them to numbers (which would strip formatting / introduce float rounding).

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.775.03'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.495.88'
$ws.Range('E3').Value = '  -0.55%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '602.20'
$ws.Range('E5').Value = '  -1.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.46'
$ws.Range('E6').Value = '  -3.35%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.495.02'
$ws.Range('E7').Value = '  -0.52%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').Value = '  -1.81%  '
$ws.Range('E10').Value = '  -1.23%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.83'
$ws.Range('E11').Value = '  +2.73%  '
$ws.Range('E12').Value = '  -2.21%  '
$ws.Range('E13').Value = '  -1.96%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.090.72'
$ws.Range('E14').Value = '  -0.44%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '31.21'
$ws.Range('E15').Value = '  -4.62%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.494.19'
$ws.Range('E16').Value = '  -0.69%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '66.778.07'
$ws.Range('E17').Value = '  -1.25%  '
$ws.Range('E18').Value = '  -0.81%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.52'
$ws.Range('E19').Value = '  +6.42%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.38'
$ws.Range('E20').Value = '  -3.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.32'
$ws.Range('E21').Value = '  -1.78%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '433.44'
$ws.Range('E22').Value = '  -3.49%  '
$ws.Range('E23').Value = '  -4.39%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '79.76'
$ws.Range('E24').Value = '  +1.89%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.633.81'
$ws.Range('E25').Value = '  -0.55%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('E27').Value = '  -3.03%  '
$ws.Range('E28').Value = '  -7.57%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.79'
$ws.Range('E29').Value = '  -3.37%  '
$ws.Range('E30').Value = '  -7.47%  '
$ws.Range('E31').Value = '  -1.11%  '
$ws.Range('E32').Value = '  -4.06%  '
$ws.Range('E33').Value = '  -0.05%  '
$ws.Range('E34').Value = '  -2.44%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '25.31'
$ws.Range('E35').Value = '  -1.90%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.490.69'
$ws.Range('E36').Value = '  -0.50%  '
$ws.Range('E37').Value = '  -4.75%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.86'
$ws.Range('E38').Value = '  -5.55%  '
$ws.Range('E39').Value = '  -0.89%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0889'
$ws.Range('E42').Value = '  -1.10%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '168.81'
$ws.Range('E43').Value = '  -3.11%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.09'
$ws.Range('E44').Value = '  -9.35%  '
$ws.Range('E45').Value = '  -1.63%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.896'
$ws.Range('E46').Value = '  +1.38%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '28.97'
$ws.Range('E47').Value = '  -4.61%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '45.70'
$ws.Range('E48').Value = '  -2.22%  '
$ws.Range('E49').Value = '  +0.70%  '
$ws.Range('E50').Value = '  -3.05%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.41'
$ws.Range('E51').Value = '  -4.95%  '
